# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-14 03:24:33
#
# For every cell in the "Recorded By" column (G) whose comma-separated
# value begins with "System" (or "system"), move all "System"/"system"
# tokens to the end of the list, preserving the relative order of the
# remaining tokens (and of the moved tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$startRow = $used.Row
$rowCount = $used.Rows.Count
$lastRow = $startRow + $rowCount - 1

for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val.GetType().Name -ne "String") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")
    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    # Only touch rows whose first token is literally "System"/"system".
    if ($trimmedParts[0].ToLower() -ne "system") { continue }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmedParts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $final = $otherParts + $systemParts
    $newVal = [string]::Join(", ", $final)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
